$wb = $excel.ActiveWorkbook

# Update "展览" (Exhibition) sheet, row 3: F3 78 -> 79, G3 29.9 -> 50
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 79
$ws1.Range("G3").Value = 50

# Update "全部类型" (All Types) sheet, row 3: F3 78 -> 79, G3 29.9 -> 50
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 79
$ws4.Range("G3").Value = 50
